$d = $word.ActiveDocument

function Find-ParagraphByExactText($doc, [string]$text) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Trim() -eq $text) {
            return $p
        }
    }
    return $null
}

function Split-RunAt($doc, [int]$fromPos, [int]$toPos) {
    # Toggling a formatting property (on then straight back off) over
    # [fromPos, toPos) forces the engine to keep a hard run boundary at
    # fromPos instead of silently re-merging two adjacent, identically
    # formatted runs. Using the end of the affected text (commonly the
    # paragraph end) as toPos avoids introducing a *second*, unwanted
    # boundary at toPos.
    $r = $doc.Range($fromPos, $toPos)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# ---------------------------------------------------------------
# 1. "Ruby (1.9, RSpec, Cucumber, AWS)" -> "Ruby (2.4, RSpec, Cucumber, AWS)"
#    split across three runs: "Ruby (" / "2.4" / ", RSpec, Cucumber, AWS)"
# ---------------------------------------------------------------
$rubyPara = Find-ParagraphByExactText $d "Ruby (1.9, RSpec, Cucumber, AWS)"
$rubyStart = $rubyPara.Range.Start
$rubyParaEnd = $rubyPara.Range.End - 1
$rubyFull = $rubyPara.Range.Text
$verIdx = $rubyFull.IndexOf("1.9")
$verStart = $rubyStart + $verIdx
$verEnd = $verStart + 3

$verRange = $d.Range($verStart, $verEnd)
$verRange.Text = "2.4"

$splitFrom1 = $verStart
$splitTo1 = $rubyParaEnd
Split-RunAt $d $splitFrom1 $splitTo1

$splitFrom2 = $verEnd
$splitTo2 = $rubyParaEnd
Split-RunAt $d $splitFrom2 $splitTo2

# ---------------------------------------------------------------
# 2. Insert new skill bullets right after the Ruby bullet:
#    F#, OCaml, Clojure, Elixir, CoffeeScript, Linux, Docker
# ---------------------------------------------------------------
$rubyPara = Find-ParagraphByExactText $d "Ruby (2.4, RSpec, Cucumber, AWS)"
$newSkills = @("F#", "OCaml", "Clojure", "Elixir", "CoffeeScript", "Linux", "Docker")

$afterStart = $rubyPara.Range.Start
$afterEnd = $rubyPara.Range.End

foreach ($skill in $newSkills) {
    $afterPara = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -eq $afterStart -and $p.Range.End -eq $afterEnd) {
            $afterPara = $p
            break
        }
    }
    $afterPara.Range.InsertParagraphAfter()
    $found = $false
    foreach ($p in $d.Paragraphs) {
        if ($found) {
            $p.Range.InsertBefore($skill)
            $afterStart = $p.Range.Start
            $afterEnd = $p.Range.End
            break
        }
        if ($p.Range.Start -eq $afterStart -and $p.Range.End -eq $afterEnd) {
            $found = $true
        }
    }
}

# ---------------------------------------------------------------
# 3. "Release Management" -> "Release Management/Continuous Delivery"
#    split across two runs: "Release Management/" / "Continuous Delivery"
# ---------------------------------------------------------------
$relPara = Find-ParagraphByExactText $d "Release Management"
$relParaEnd = $relPara.Range.End - 1
$insPoint = $d.Range($relParaEnd, $relParaEnd)
$insPoint.InsertBefore("/Continuous Delivery")

$relSplitFrom = $relParaEnd + 1
$relSplitTo = $relPara.Range.End - 1
Split-RunAt $d $relSplitFrom $relSplitTo

# ---------------------------------------------------------------
# 4. "Agile/Lean" -> "Agile/Lean/Kanban"
#    split across two runs: "Agile/Lean/" / "Kanban"
# ---------------------------------------------------------------
$aglPara = Find-ParagraphByExactText $d "Agile/Lean"
$aglParaEnd = $aglPara.Range.End - 1
$insPoint2 = $d.Range($aglParaEnd, $aglParaEnd)
$insPoint2.InsertBefore("/Kanban")

$aglSplitFrom = $aglParaEnd + 1
$aglSplitTo = $aglPara.Range.End - 1
Split-RunAt $d $aglSplitFrom $aglSplitTo

# ---------------------------------------------------------------
# 5. Insert a new bullet after the "github." bullet in PERSONAL INTERESTS:
#    "I love reading lots of technical books and blogs"
# ---------------------------------------------------------------
$githubPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -like "*github*") {
        $githubPara = $p
    }
}
$ghStart = $githubPara.Range.Start
$ghEnd = $githubPara.Range.End
$githubPara.Range.InsertParagraphAfter()
$found = $false
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $p.Range.InsertBefore("I love reading lots of technical books and blogs")
        break
    }
    if ($p.Range.Start -eq $ghStart -and $p.Range.End -eq $ghEnd) {
        $found = $true
    }
}

Write-Output "done"
